# Capstone Hour Tracker - "a lot, check hour tracker"
#
# The author filled in a new work-log entry. What used to be an empty
# template row (row 18) now holds a real start/end time plus the
# "what I accomplished" / "what's next" notes, and a short aside note
# was added next to the existing row 17 entry (column J). The running
# total / goal formulas in row 4 and row 7 already reference these
# cells, so they recalculate automatically once the new values land.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New aside note on the existing row 17 entry -------------------------
# (set first so its shared-string entry lands before the row-18 strings,
# matching the order new strings were appended in the real edit)
$ws.Range("J17").Value2 = "It's weird to think this is a shared document. Not through excel, but github lol"

# --- Fill in row 18 (previously a blank template row) --------------------
$ws.Range("B18").Value2 = 45210                  # Date: 10/11/2023
$ws.Range("C18").Value2 = 0.77569444444444446    # Start Time: 6:37 PM
$ws.Range("D18").Value2 = 0.86458333333333337    # End Time: 8:45 PM

# D18 picked up the plain "General" column-default style; copy the time
# number format (and matching border/fill) from D17 so it renders as a
# time like the rest of column D.
$ws.Range("D17").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("G18").Value2 = "Finished up the view panels, started on the sidebar, started laying out the editor page"
$ws.Range("H18").Value2 = "Next items on the horizon are: file heirarchy viewer, theme switcher, settings, terminal, top bar, moving most state into a context, as well as starting on a lot of backend type shit"

# --- Row heights follow the new/changed wrapped content -------------------
$ws.Rows.Item(1).RowHeight = 36
$ws.Rows.Item(3).RowHeight = 30
$ws.Rows.Item(5).RowHeight = 60
$ws.Rows.Item(6).RowHeight = 75
$ws.Rows.Item(7).RowHeight = 30
$ws.Rows.Item(8).RowHeight = 60
$ws.Rows.Item(9).RowHeight = 60
$ws.Rows.Item(10).RowHeight = 30
$ws.Rows.Item(11).RowHeight = 90
$ws.Rows.Item(12).RowHeight = 60
$ws.Rows.Item(13).RowHeight = 30
$ws.Rows.Item(14).RowHeight = 90
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(16).RowHeight = 30
$ws.Rows.Item(17).RowHeight = 90
$ws.Rows.Item(18).RowHeight = 75

# --- Selection ends up on D19 after entering the new row's data ----------
$ws.Range("D19").Select()
